$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MW CNN 1D")
$ws.Activate()

# New "250 Packets" summary block header (N1) and Acc/Loss/Time sub-headers (M2:O2)
$ws.Range("N1").Value = "250 Packets"
$ws.Range("M2").Value = "Acc"
$ws.Range("N2").Value = "Loss"
$ws.Range("O2").Value = "Time"

# Fill in the 270-packet run data (I3:K51), previously only headered with no rows
$rows = @(
    @(88.526570796966496,0.27364717987528597,571.54140663146904),
    @(89.508384466171194,0.26449414726576298,563.20081591606095),
    @(88.806414604186998,0.25839889372028602,570.53296017646699),
    @(89.207202196121202,0.26261496832917403,567.04648327827397),
    @(89.029335975646902,0.248651185612245,567.60861730575505),
    @(88.704437017440796,0.26690111964305602,567.39737844467095),
    @(88.981902599334703,0.24983153020603099,569.47559022903397),
    @(89.453840255737305,0.25886353362694098,567.98920059204102),
    @(89.5131289958953,0.27139978508277002,566.18249654769897),
    @(89.226174354553194,0.25948665560720902,568.16771006584099),
    @(88.661748170852604,0.26117336892303,566.52932882308903),
    @(88.910758495330796,0.26289581964776798,567.09138488769497),
    @(88.161361217498694,0.28794452782073499,568.954289674758),
    @(89.024591445922795,0.27293423157291502,564.79920411109902),
    @(89.811939001083303,0.25717963969796198,570.29118657112099),
    @(88.801670074462805,0.28394434498728299,564.91241025924603),
    @(88.294166326522799,0.28930436779771002,565.87722325325001),
    @(89.088624715804997,0.26881124585344601,568.35680985450699),
    @(88.391393423080402,0.291519165259574,569.515851259231),
    @(87.938433885574298,0.28873213954559801,569.02983236312798),
    @(88.763725757598806,0.26834346629382499,566.78986501693703),
    @(88.009577989578204,0.274134513162546,568.81757926940895),
    @(89.155024290084796,0.259638950536135,568.20859313011101),
    @(88.894158601760793,0.265488644609297,569.044016361236),
    @(89.444351196289006,0.24944659968504601,567.303020238876),
    @(88.704437017440796,0.26888305771927401,564.03714585304203),
    @(87.675195932388306,0.27972395518576798,570.32110118865899),
    @(90.011143684387207,0.23787439667117899,565.27814984321503),
    @(88.730525970458899,0.27651749239611101,569.52728033065796),
    @(82.941639423370304,2.6234215325207999,567.20298838615395),
    @(88.590604066848698,0.25385649970335,568.78105330467201),
    @(89.109969139099107,0.26093641278191998,564.85970544814995),
    @(88.502860069274902,0.279851582453916,565.68032455444302),
    @(87.888634204864502,0.276040247304353,564.46240520477295),
    @(88.939219713211003,0.26452549647477203,566.59601569175697),
    @(89.9874329566955,0.25380408502098301,568.67201066017105),
    @(89.280718564987097,0.25534665625980801,568.89972043037403),
    @(87.962150573730398,0.28369151539454701,566.00424218177795),
    @(89.088624715804997,0.25009892779530601,568.02661538124005),
    @(89.079135656356797,0.26349528033048503,565.79061031341496),
    @(88.080728054046602,0.28729819680042201,570.74440741539001),
    @(88.474398851394596,0.27594648872458599,566.73435997962895),
    @(89.107596874237004,0.25494215230629202,567.932703971862),
    @(89.164513349532996,0.26137452160043101,570.83284282684303),
    @(89.254629611968994,0.25538730008993299,565.16936302185002),
    @(88.991391658782902,0.29412961411902799,569.66812443733204),
    @(89.041191339492798,0.26776675705595399,567.20592927932705),
    @(89.155024290084796,0.249937309914257,567.783757925033),
    @(88.611948490142794,0.28223201247204599,570.62241840362503)
)

$arr = New-Object 'object[,]' $rows.Count,3
for ($i = 0; $i -lt $rows.Count; $i++) {
    for ($j = 0; $j -lt 3; $j++) {
        $arr[$i,$j] = $rows[$i][$j]
    }
}
$ws.Range("I3:K51").Value = $arr

$ws.Range("S16").Select() | Out-Null
